# Add a new column R (year 2021) to the unemployment-rate table on sheet1.
# Mirrors the existing year columns (D:Q, years 2007-2020) by copying
# formatting from a same-styled neighbour cell in each row, then writing
# the 2021 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Header row (year label 2021), formatted like the other year headers ---
$ws.Cells.Item(4, 18).Value = 2021
[void]$ws.Cells.Item(4, 17).Copy()
[void]$ws.Cells.Item(4, 18).PasteSpecial($xlPasteFormats)

# --- Data rows: row number -> value (or $null for a blank cell) ---
# "Bold" rows (country/oblast totals) copy formatting from A5 (bold, no
# number-format style already present in the sheet); "plain" rows copy
# formatting from A10 (regular, no number-format style).
$values = @{
    5  = 5.3
    6  = 6.3
    7  = 4.7
    8  = $null
    9  = 6.6
    10 = 7.5
    11 = 6.2
    12 = 11.8
    13 = 15.5
    14 = 9.6999999999999993
    15 = 6.3
    16 = 7.5
    17 = 5.6
    18 = 6.3
    19 = 10.8
    20 = 4.3
    21 = 1.9
    22 = 3.1
    23 = 1.1000000000000001
    24 = 2.6
    25 = 3.8
    26 = 1.7
    27 = 5.3
    28 = 6.2
    29 = 4.8
    30 = 4.0999999999999996
    31 = 3.3
    32 = 4.9000000000000004
    33 = 2.8
    34 = 3.4
    35 = 2.6
    36 = $null
    37 = 15.7
    38 = 7.9
    39 = 4.5
    40 = 4.4000000000000004
    41 = 2.9
    42 = 1.4
}

$boldRows = @(5, 9, 12, 15, 18, 21, 24, 27, 30, 33)

for ($row = 5; $row -le 42; $row++) {

    if ($boldRows -contains $row) {
        [void]$ws.Cells.Item(5, 1).Copy()
    }
    else {
        [void]$ws.Cells.Item(10, 1).Copy()
    }
    [void]$ws.Cells.Item($row, 18).PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    $v = $values[$row]
    if ($null -ne $v) {
        $ws.Cells.Item($row, 18).Value = $v
    }
}

# --- Last row (70 and over): data not available, shows an ellipsis with
#     the bottom-border style used by the rest of that row ---
$ws.Cells.Item(43, 18).Value = "…"
[void]$ws.Cells.Item(43, 1).Copy()
[void]$ws.Cells.Item(43, 18).PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Cells.Item(43, 18).HorizontalAlignment = -4152

# --- Move the active selection like in the saved workbook ---
[void]$ws.Range("S1").Select()
